$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 340.35294
$ws.Range("I33").Value = 351.77777
$ws.Range("K33").Value = 351.77777
$ws.Range("M33").Value = -122.77777
# Row 63
$ws.Range("H63").Value = 10246
$ws.Range("I63").Value = 10246
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 10246
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -9622
$ws.Range("N63").ClearContents()
# Row 66
$ws.Range("H66").Value = 10246
$ws.Range("I66").Value = 10246
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 30738
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -27618
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10472.845
$ws.Range("I32").Value = 8181.4653
$ws.Range("J32").Value = 17041.467
$ws.Range("K32").Value = 8181.4653
$ws.Range("L32").Value = 17041.467
$ws.Range("M32").Value = -7894.4653
$ws.Range("N32").Value = -17615.467
# Row 45
$ws.Range("H45").Value = 2201.7144
$ws.Range("I45").Value = 2068.6667
$ws.Range("K45").Value = 2068.6667
$ws.Range("M45").Value = -1691.6667
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 74
$ws.Range("H74").Value = 6301.9355
$ws.Range("I74").Value = 3393.5833
$ws.Range("K74").Value = 3393.5833
$ws.Range("M74").Value = -2519.5833
# Row 77
$ws.Range("H77").Value = 6301.9355
$ws.Range("I77").Value = 3393.5833
$ws.Range("K77").Value = 16967.9165
$ws.Range("M77").Value = -12599.9165
# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
# Row 130
$ws.Range("H130").Value = 59966.668
$ws.Range("J130").Value = 59966.668
$ws.Range("L130").Value = 59966.668
$ws.Range("N130").Value = -70006.66800000001
# Row 132
$ws.Range("H132").Value = 5522.05
$ws.Range("I132").Value = 1262.0588
$ws.Range("J132").Value = 8670.739
$ws.Range("K132").Value = 3786.1764
$ws.Range("L132").Value = 26012.217
$ws.Range("M132").Value = -1256.1764
$ws.Range("N132").Value = -31072.217

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 522.2222
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = -327
$ws.Range("N22").Value = -946
# Row 54
$ws.Range("H54").Value = 15986.637
$ws.Range("I54").Value = 3941.5
$ws.Range("J54").Value = 22869.572
$ws.Range("K54").Value = 3941.5
$ws.Range("L54").Value = 22869.572
$ws.Range("M54").Value = -3457.5
$ws.Range("N54").Value = -23837.572

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 959.1429000000001
$ws.Range("I16").Value = 750.25
$ws.Range("J16").Value = 1237.6666
$ws.Range("K16").Value = 750.25
$ws.Range("L16").Value = 1237.6666
$ws.Range("M16").Value = -463.25
$ws.Range("N16").Value = -1811.6666
# Row 45
$ws.Range("H45").Value = 13333
$ws.Range("I45").Value = 9999
$ws.Range("K45").Value = 9999
$ws.Range("M45").Value = -9406
# Row 58
$ws.Range("H58").Value = 1979369.9
$ws.Range("I58").Value = 3137131.5
$ws.Range("J58").Value = 4364.647
$ws.Range("K58").Value = 3137131.5
$ws.Range("L58").Value = 4364.647
$ws.Range("M58").Value = -3136928.5
$ws.Range("N58").Value = -4770.647
# Row 107
$ws.Range("H107").Value = 571.5
$ws.Range("J107").Value = 678.25
$ws.Range("L107").Value = 678.25
$ws.Range("N107").Value = -4518.25
# Row 113
$ws.Range("H113").Value = 959.1429000000001
$ws.Range("I113").Value = 750.25
$ws.Range("J113").Value = 1237.6666
$ws.Range("K113").Value = 750.25
$ws.Range("L113").Value = 1237.6666
$ws.Range("M113").Value = 1419.75
$ws.Range("N113").Value = -5577.6666
# Row 134
$ws.Range("H134").Value = 4577.129
$ws.Range("I134").Value = 3710.111
$ws.Range("J134").Value = 4931.8184
$ws.Range("K134").Value = 11130.333
$ws.Range("L134").Value = 14795.4552
$ws.Range("M134").Value = -8595.332999999999
$ws.Range("N134").Value = -19865.4552
# Row 136
$ws.Range("H136").Value = 1979369.9
$ws.Range("I136").Value = 3137131.5
$ws.Range("J136").Value = 4364.647
$ws.Range("K136").Value = 9411394.5
$ws.Range("L136").Value = 13093.941
$ws.Range("M136").Value = -9408844.5
$ws.Range("N136").Value = -18193.941

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 34800.535
$ws.Range("I131").Value = 758.7857
$ws.Range("J131").Value = 68842.28999999999
$ws.Range("K131").Value = 2276.3571
$ws.Range("L131").Value = 206526.87
$ws.Range("M131").Value = 2763.6429
$ws.Range("N131").Value = -216606.87
# Row 132
$ws.Range("H132").Value = 1682.963
$ws.Range("I132").Value = 1455.6923
$ws.Range("J132").Value = 1894
$ws.Range("K132").Value = 13101.2307
$ws.Range("L132").Value = 17046
$ws.Range("M132").Value = -10571.2307
$ws.Range("N132").Value = -22106
# Row 133
$ws.Range("H133").Value = 3575.182
$ws.Range("I133").Value = 3258.5557
$ws.Range("J133").Value = 5000
$ws.Range("K133").Value = 9775.667099999999
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = -4715.667099999999
$ws.Range("N133").Value = -25120
# Row 134
$ws.Range("H134").Value = 6046
$ws.Range("I134").Value = 6046
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 18138
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -13068
$ws.Range("N134").ClearContents()
# Row 137
$ws.Range("H137").Value = 33370688
$ws.Range("I137").Value = 45456388
$ws.Range("J137").Value = 135008
$ws.Range("K137").Value = 136369164
$ws.Range("L137").Value = 405024
$ws.Range("M137").Value = -136364064
$ws.Range("N137").Value = -415224
# Row 138
$ws.Range("H138").Value = 7612.647
$ws.Range("I138").Value = 11072
$ws.Range("J138").Value = 2670.7144
$ws.Range("K138").Value = 33216
$ws.Range("L138").Value = 8012.1432
$ws.Range("M138").Value = -28076
$ws.Range("N138").Value = -18292.1432

$ws = $wb.Worksheets.Item("GSM")
# Row 64
$ws.Range("H64").Value = 37270.5
$ws.Range("J64").Value = 37270.5
$ws.Range("L64").Value = 37270.5
$ws.Range("N64").Value = -37766.5
# Row 67
$ws.Range("H67").Value = 37270.5
$ws.Range("J67").Value = 37270.5
$ws.Range("L67").Value = 37270.5
$ws.Range("N67").Value = -38986.5
# Row 70
$ws.Range("H70").Value = 4867.6333
$ws.Range("I70").Value = 4187.875
$ws.Range("J70").Value = 5114.8184
$ws.Range("K70").Value = 4187.875
$ws.Range("L70").Value = 5114.8184
$ws.Range("M70").Value = -3917.875
$ws.Range("N70").Value = -5654.8184
# Row 73
$ws.Range("H73").Value = 4867.6333
$ws.Range("I73").Value = 4187.875
$ws.Range("J73").Value = 5114.8184
$ws.Range("K73").Value = 4187.875
$ws.Range("L73").Value = 5114.8184
$ws.Range("M73").Value = -3251.875
$ws.Range("N73").Value = -6986.8184
# Row 126
$ws.Range("H126").Value = 3019.95
$ws.Range("I126").Value = 1799.9166
$ws.Range("J126").Value = 4850
$ws.Range("K126").Value = 5399.7498
$ws.Range("L126").Value = 14550
$ws.Range("M126").Value = -2929.7498
$ws.Range("N126").Value = -19490
# Row 132
$ws.Range("H132").Value = 6803.8096
$ws.Range("I132").Value = 9727.166999999999
$ws.Range("J132").Value = 2906
$ws.Range("K132").Value = 29181.501
$ws.Range("L132").Value = 8718
$ws.Range("M132").Value = -26651.501
$ws.Range("N132").Value = -13778

$ws = $wb.Worksheets.Item("LTW")
# Row 94
$ws.Range("H94").Value = 20000
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("N94").Value = -21352

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1392
$ws.Range("I126").Value = 1063
$ws.Range("J126").Value = 1830.6666
$ws.Range("K126").Value = 3189
$ws.Range("L126").Value = 5491.9998
$ws.Range("M126").Value = -719
$ws.Range("N126").Value = -10431.9998
# Row 132
$ws.Range("H132").Value = 1478.7428
$ws.Range("I132").Value = 1346.5186
$ws.Range("J132").Value = 1925
$ws.Range("K132").Value = 4039.5558
$ws.Range("L132").Value = 5775
$ws.Range("M132").Value = -1509.5558
$ws.Range("N132").Value = -10835
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
